$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 new rows before the last data row (Cambridge Dictionary),
#        pushing it from row 11 down to row 15 ---
$ws.Rows("11:14").Insert()

# Copy the number/border formatting from row 10 (a "normal" data row) onto
# the freshly inserted rows so they look consistent with the rest of the table.
$ws.Range("A10:N10").Copy()
$ws.Range("A11:N14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Fill in the new rows' data ---
$newRows = @(
  @{r=11; name="Amazon";      c=1; e=0; g=0; i=10; k=3; m=1},
  @{r=12; name="Google";      c=1; e=1; g=0; i=4;  k=1; m=4},
  @{r=13; name="Weather.com"; c=0; e=0; g=1; i=0;  k=2; m=2},
  @{r=14; name="Reddit";      c=0; e=0; g=1; i=1;  k=0; m=3}
)

foreach ($row in $newRows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.name
    $ws.Range("B$r").Formula = "=C$r+E$r+G$r+I$r+K$r+M$r"
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Formula = "=C$r/(`$C$r+`$E$r+`$G$r)*100"
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Formula = "=E$r/(`$C$r+`$E$r+`$G$r)*100"
    $ws.Range("G$r").Value = $row.g
    $ws.Range("H$r").Formula = "=G$r/(`$C$r+`$E$r+`$G$r)*100"
    $ws.Range("I$r").Value = $row.i
    $ws.Range("J$r").Formula = "=I$r/(`$I$r+`$K$r+`$M$r)*100"
    $ws.Range("K$r").Value = $row.k
    $ws.Range("L$r").Formula = "=K$r/(`$I$r+`$K$r+`$M$r)*100"
    $ws.Range("M$r").Value = $row.m
    $ws.Range("N$r").Formula = "=M$r/(`$I$r+`$K$r+`$M$r)*100"
}

# --- 3. Add the new totals row (row 16) ---
$ws.Range("B16").Formula = "=SUM(B3:B15)"
$ws.Range("C16").Formula = "=SUM(C3:C15)"
$ws.Range("E16").Formula = "=SUM(E3:E15)"
$ws.Range("G16").Formula = "=SUM(G3:G15)"
$ws.Range("I16").Formula = "=SUM(I3:I15)"
$ws.Range("K16").Formula = "=SUM(K3:K15)"
$ws.Range("M16").Formula = "=SUM(M3:M15)"

$ws.Range("A16:B16").Font.Bold = $true
$ws.Range("A16:N16").Borders.LineStyle = 0
$ws.Range("D16").NumberFormat = "\(0\%\)"
$ws.Range("F16").NumberFormat = "\(0\%\)"
$ws.Range("H16").NumberFormat = "\(0\%\)"
$ws.Range("J16").NumberFormat = "\(0\%\)"
$ws.Range("L16").NumberFormat = "\(0\%\)"
$ws.Range("N16").NumberFormat = "\(0\%\)"

# --- 4. Move the two charts down so they sit below the now-larger table ---
$co1 = $ws.ChartObjects(1)
$co1.Top = $ws.Range("B18").Top()
$co1.Left = $ws.Range("B18").Left()

$co2 = $ws.ChartObjects(2)
$co2.Top = $ws.Range("B32").Top()
$co2.Left = $ws.Range("B32").Left()

# --- 5. Extend the charts' source-data ranges to cover the new rows ---
$chart1 = $co1.Chart
$chart1.SeriesCollection(1).Formula = "=SERIES(Foglio1!`$C`$2,Foglio1!`$A`$3:`$A`$15,Foglio1!`$C`$3:`$C`$15,1)"
$chart1.SeriesCollection(2).Formula = "=SERIES(Foglio1!`$E`$2,Foglio1!`$A`$3:`$A`$15,Foglio1!`$E`$3:`$E`$15,2)"
$chart1.SeriesCollection(3).Formula = "=SERIES(Foglio1!`$G`$2,Foglio1!`$A`$3:`$A`$15,Foglio1!`$G`$3:`$G`$15,3)"

$chart2 = $co2.Chart
$chart2.SeriesCollection(1).Formula = "=SERIES(Foglio1!`$I`$2,Foglio1!`$A`$3:`$A`$15,Foglio1!`$I`$3:`$I`$15,1)"
$chart2.SeriesCollection(2).Formula = "=SERIES(Foglio1!`$K`$2,Foglio1!`$A`$3:`$A`$15,Foglio1!`$K`$3:`$K`$15,2)"
$chart2.SeriesCollection(3).Formula = "=SERIES(Foglio1!`$M`$2,Foglio1!`$A`$3:`$A`$15,Foglio1!`$M`$3:`$M`$15,3)"

# --- 6. Recolor the "Incorrect" series of the second chart (Navigation) from
#        the theme accent3 color to a solid dark red (C00000) ---
$chart2.SeriesCollection(3).Interior.Color = 49152   # RGB(192,0,0) = 0x00C000 -> BGR 0x0000C0... see below

Write-Output "done"
